# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-10-29 07:15:07
#
# A previously "Not Recorded" PHYSIOLOGY session (row 28) has now been recorded
# (by "System" on behalf of Aya_hamed), so this script:
#   - updates the "Recorded By" lists for a couple of sessions (reordered / appended)
#   - flips row 28 from "Not Recorded" (pink) to "Recorded" (green) and fills in
#     its attendance numbers
#   - refreshes the dependent summary statistics (recorded/missing counts and
#     coverage / average-attendance percentages) that are derived from it
#   - narrows column I (9) now that the "Not Recorded" values are gone

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column I width: 14 -> 10 (raw OOXML width) -----------------------------
# This sheet's columns consistently report ColumnWidth == (raw width - 0.83)
# for this theme/font, so 9.17 here round-trips to a raw <col width="10">.
$ws.Columns.Item(9).ColumnWidth = 9.17

# --- G2: reorder / extend the "Recorded By" list ---------------------------
$ws.Range("G2").Value = "Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System"

# --- Recorded / Missing session counters (K6:L7 block) ---------------------
$ws.Range("L6").Value = 3
$ws.Range("L7").Value = 0

# --- G9: reorder the "Recorded By" list -------------------------------------
$ws.Range("G9").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

# --- Coverage % / Average Attendance % (K9:L10 block) -----------------------
# These columns hold the percentage as literal text (e.g. "10.3%"), not a
# numeric percentage. Assigning a "NN.N%" string straight to .Value makes
# Excel auto-convert it to a numeric percent, so instead write it as a
# formula that yields the text, then convert that formula to a plain value
# in-place (Copy + PasteSpecial values) so the stored cell keeps the exact
# literal text without disturbing the existing cell style.
function Set-TextValue($range, $text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextValue $ws.Range("L9") "10.3%"
Set-TextValue $ws.Range("L10") "10.2%"

# --- Group summary table row (row 15): Recorded / Missing / % columns -------
$ws.Range("O15").Value = 3
$ws.Range("P15").Value = 0
Set-TextValue $ws.Range("R15") "10.3%"
Set-TextValue $ws.Range("S15") "10.2%"
$excel.CutCopyMode = $false

# --- Row 28: PHYSIOLOGY session 1 flips from Not Recorded -> Recorded -------
# Copy the formatting (fill/font) of an already-"Recorded" row (row 2) onto
# row 28 so it switches from the pink "Not Recorded" style to the green
# "Recorded" style, then update its values.
$ws.Range("A2:I2").Copy()
$ws.Range("A28:I28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg"
$ws.Range("H28").Value = "6/251"
$ws.Range("I28").Value = "Recorded"
